# "Generate Report for Handback" - update the localization-status report to
# reflect that the de-de / zh-cn handback cycle has completed: the overview
# status flips from "Ready for handoff" to "Handed back: in sync with en-US",
# and each language sheet gets its Latest Target File / Latest Handback File
# / Latest Handback DateTime columns populated (with a real hyperlink on the
# target-file cell), plus the resulting column widths widen to fit the new
# text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/23a5b0b50f176275524b388c5534f8a54b32c36e/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both file rows
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status column (C) picks up the same new text
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# Latest Target File (I) now links to a.md, like the Source File Name column
$zhcn.Range("I2").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $baseUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$zhcn.Range("I3").Value = "a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $baseUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null

# Latest Handback File (J) + Latest Handback DateTime (K)
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-18 00:34:32"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-18 00:34:32"

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("I2"), $baseUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$dede.Range("I3").Value = "a.md"
$dede.Hyperlinks.Add($dede.Range("I3"), $baseUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null

$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-18 00:34:39"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-18 00:34:39"

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(10).ColumnWidth = 39.17
